$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.873.13"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "2.413.75"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'551.50"
$ws.Range("E5").Value = "  -0.52%  "
$ws.Range("D6").Value = "'137.17"
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.590"
$ws.Range("E8").Value = "  +3.67%  "
$ws.Range("E9").Value = "  -2.02%  "
$ws.Range("D10").Value = "'5.69"
$ws.Range("E10").Value = "  -1.91%  "
$ws.Range("E11").Value = "  -1.03%  "
$ws.Range("E12").Value = "  -2.39%  "
$ws.Range("D13").Value = "'25.59"
$ws.Range("E13").Value = "  +3.89%  "
$ws.Range("D14").Value = "2.845.83"
$ws.Range("D15").Value = "59.824.46"
$ws.Range("E15").Value = "  +0.43%  "
$ws.Range("E16").Value = "  -1.50%  "
$ws.Range("D17").Value = "2.433.23"
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("E19").Value = "  -1.09%  "
$ws.Range("D20").Value = "'328.89"
$ws.Range("E20").Value = "  -2.08%  "
$ws.Range("D21").Value = "'6.65"
$ws.Range("E21").Value = "  -4.62%  "
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").Value = "'66.50"
$ws.Range("E23").Value = "  +2.95%  "
$ws.Range("D24").Value = "'0.171"
$ws.Range("E24").Value = "  +1.12%  "
$ws.Range("D25").Value = "'8.64"
$ws.Range("E25").Value = "  +1.56%  "
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("D28").Value = "0.0₃0774"
$ws.Range("E28").Value = "  -1.03%  "
$ws.Range("E29").Value = "  -1.89%  "
$ws.Range("D30").Value = "'168.12"
$ws.Range("E30").Value = "  -1.43%  "
$ws.Range("D31").Value = "'6.10"
$ws.Range("E31").Value = "  -2.33%  "
$ws.Range("E32").Value = "  -0.49%  "
$ws.Range("E33").Value = "  -1.18%  "
$ws.Range("E35").Value = "  -0.79%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "'4.22"
$ws.Range("E37").Value = "  -2.11%  "
$ws.Range("E38").Value = "  -2.66%  "
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").Value = "'313.56"
$ws.Range("E39").Value = "  +2.75%  "
$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D40").Value = "'0.408"
$ws.Range("E40").Value = "  -2.70%  "
$ws.Range("E41").Value = "  -1.98%  "
$ws.Range("D42").Value = "'138.62"
$ws.Range("E42").Value = "  -2.66%  "
$ws.Range("E43").Value = "  +0.39%  "
$ws.Range("D44").Value = "'0.0517"
$ws.Range("E44").Value = "  -1.82%  "
$ws.Range("D45").Value = "'19.50"
$ws.Range("E45").Value = "  +2.17%  "
$ws.Range("E46").Value = "  +1.53%  "
$ws.Range("E47").Value = "  -0.95%  "
$ws.Range("D48").Value = "'0.386"
$ws.Range("E48").Value = "  -4.88%  "
$ws.Range("D49").Value = "'17.67"
$ws.Range("E49").Value = "  -1.45%  "
$ws.Range("D50").Value = "'11.05"
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("D51").Value = "'1.57"
$ws.Range("E51").Value = "  -1.78%  "
